$wb = $excel.ActiveWorkbook

# The event was cancelled: update the name with a "(cancelled)" suffix and
# mark the ticket price column as unavailable for sale on both the
# "展览" and "全部类型" sheets (row 2 = 丽水·LZ栗子动漫游戏嘉年华).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("C2").Value = "丽水·LZ栗子动漫游戏嘉年华（取消）"
    $ws.Range("G2").Value = "不可售"
}
